$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A49").Value = "Raffaele Prezzi"
$ws.Range("B49").Value = "Thomas Debiasi | Mai una gioia"
$ws.Range("C49").Value = "Alessio Bragagna | SHARK ATTACK"
$ws.Range("D49").Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Range("E49").Value = "Andreas Galli | SdrumALA"
$ws.Range("F49").Value = "Mattia Tezzele | U.SGUARNA"
